$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "*_old" -> "*_FV2210", "*_new" -> "*_FV2304"
# ---------------------------------------------------------------------------
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}
# column 11 ("diff") keeps its name
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table (ListObject), keeping the
#    original header-row formatting intact (bold / shaded / bordered) and
#    without Excel bolting on a header dxf or a named table style.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")

# Stash a copy of the header formatting on a scratch cell far outside the
# used range so we can restore it verbatim after the table is created.
$scratch = $ws.Range("A200")
$ws.Range("A1").Copy($scratch)

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U93")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Restore the original header formatting in one shot (so it collapses back
# onto the pre-existing style record instead of creating new ones).
$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
